$d = $word.ActiveDocument

# Update the date paragraph (first paragraph in the document, outside the table)
$d.Paragraphs.Item(1).Range.Find.Execute("2024-10-04 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-05 Saturday", 2) | Out-Null

# Update each table cell value directly via its Range, to avoid any ambiguity
# from duplicate values (e.g. "14-4=10" appears twice with different replacements).
$t = $d.Tables.Item(1)
$newValues = @(
    "12+23=35",
    "6+56=62",
    "34+52=86",
    "12+7=19",
    "1+75=76",
    "62-35=27",
    "3+51=54",
    "11+37=48",
    "48-32=16",
    "51-2=49",
    "61-17=44",
    "59+14=73",
    "23+44=67",
    "2+29=31",
    "62+36=98",
    "5+82=87",
    "4+13=17",
    "97-25=72",
    "47+44=91",
    "22+21=43",
    "98-2=96",
    "73-45=28",
    "70-22=48",
    "27+59=86",
    "0+87=87",
    "63+3=66",
    "87-67=20",
    "30-15=15",
    "3+29=32",
    "68-29=39",
    "19-11=8",
    "25-21=4",
    "83-30=53",
    "27+26=53",
    "97-93=4",
    "34-7=27",
    "85-65=20",
    "95-27=68",
    "0+22=22",
    "45+19=64",
    "5+85=90",
    "93-72=21",
    "95-79=16",
    "5+88=93",
    "88-37=51",
    "23+50=73",
    "75-0=75",
    "87-38=49",
    "59-14=45",
    "38+22=60",
    "38-7=31",
    "87-20=67",
    "69-17=52",
    "58-26=32",
    "50-18=32",
    "79-48=31",
    "32+40=72",
    "89-32=57",
    "56+13=69",
    "80-75=5",
    "20-9=11",
    "55+41=96",
    "46+9=55",
    "86-36=50",
    "87-66=21",
    "86-62=24",
    "85-12=73",
    "91-58=33",
    "23+8=31",
    "29+46=75",
    "14+36=50",
    "74-28=46",
    "97-47=50",
    "89-5=84",
    "59+30=89",
    "87-6=81",
    "77+9=86",
    "45+17=62",
    "96-38=58",
    "8+58=66",
    "19-14=5",
    "58-18=40",
    "44-36=8",
    "5+67=72",
    "64-52=12",
    "72-66=6",
    "11-9=2",
    "57+20=77",
    "49+6=55",
    "16+57=73",
    "60-29=31",
    "76+21=97",
    "48-1=47",
    "98-87=11",
    "0+40=40",
    "5+35=40",
    "70-27=43",
    "65+33=98",
    "31+64=95",
    "45-7=38"
)

$nRows = $t.Rows.Count
$nCols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $nRows; $r++) {
    for ($c = 1; $c -le $nCols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"
